$d = $word.ActiveDocument

# The document's single table has a trailing "Caso de uso extendido" row
# whose second cell holds nothing but the document's automatic "_GoBack"
# bookmark. That row is being removed entirely, and the "_GoBack"
# bookmark needs to live on in the (now last) empty paragraph that
# follows the table instead.

$t = $d.Tables.Item(1)

# 1. Delete the last table row ("Caso de uso extendido").
$lastRow = $t.Rows.Last
$lastRow.Delete()

# 2. Put the "_GoBack" bookmark on the trailing empty paragraph that
#    now directly follows the table (it has no runs of its own, so we
#    briefly insert a placeholder character to anchor the bookmark on,
#    then delete the character again, leaving the bookmark behind).
$cEnd = $d.Content.End
$anchor = $d.Range($cEnd - 1, $cEnd - 1)
$anchor.InsertAfter("x")

$cEnd2 = $d.Content.End
$placeholder = $d.Range($cEnd2 - 2, $cEnd2 - 1)
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder = $d.Range($cEnd2 - 2, $cEnd2 - 1)
$placeholder.Delete()
